$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F4").Value = -2
$ws.Range("F7").Value = -3
$ws.Range("F10").Value = -5
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = -6
$ws.Range("F19").Value = -6
$ws.Range("F20").Value = -8
$ws.Range("F23").Value = -3
$ws.Range("F24").Value = -4
$ws.Range("F25").Value = -7
